$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# natmiOut/OldD4/LR-pairs_lrc2p/Nid2-Col13a1.xlsx was regenerated with an
# additional target cluster ("FAPs"), expanding the sheet from 10 data rows
# (rows 2-11) to 15 data rows (rows 2-16) and refreshing every numeric metric.
# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Nid2"
$ws.Cells.Item(2, 3).Value = "Col13a1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 26.167311
$ws.Cells.Item(2, 8).Value = 78.50193300000001
$ws.Cells.Item(2, 9).Value = 0.2570455291913796
$ws.Cells.Item(2, 10).Value = 0.2645760588013328
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3015963333333334
$ws.Cells.Item(2, 14).Value = 0.9047890000000001
$ws.Cells.Item(2, 15).Value = 0.4737402899861982
$ws.Cells.Item(2, 16).Value = 0.5687660179545673
$ws.Cells.Item(2, 17).Value = 7.891965050793001
$ws.Cells.Item(2, 18).Value = 71.02768545713701
$ws.Cells.Item(2, 19).Value = 0.1217728235387799
$ws.Cells.Item(2, 20).Value = 0.1504818714105475

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Nid2"
$ws.Cells.Item(3, 3).Value = "Col13a1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 26.167311
$ws.Cells.Item(3, 8).Value = 78.50193300000001
$ws.Cells.Item(3, 9).Value = 0.2570455291913796
$ws.Cells.Item(3, 10).Value = 0.2645760588013328
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.01594066666666667
$ws.Cells.Item(3, 14).Value = 0.047822
$ws.Cells.Item(3, 15).Value = 0.02503921704145383
$ws.Cells.Item(3, 16).Value = 0.03006173650500096
$ws.Cells.Item(3, 17).Value = 0.4171243822140001
$ws.Cells.Item(3, 18).Value = 3.754119439926001
$ws.Cells.Item(3, 19).Value = 0.00643621879495831
$ws.Cells.Item(3, 20).Value = 0.007953615765217307

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Nid2"
$ws.Cells.Item(4, 3).Value = "Col13a1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 26.167311
$ws.Cells.Item(4, 8).Value = 78.50193300000001
$ws.Cells.Item(4, 9).Value = 0.2570455291913796
$ws.Cells.Item(4, 10).Value = 0.2645760588013328
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.319091
$ws.Cells.Item(4, 14).Value = 0.638182
$ws.Cells.Item(4, 15).Value = 0.5012204929723481
$ws.Cells.Item(4, 16).Value = 0.4011722455404317
$ws.Cells.Item(4, 17).Value = 8.349753434301
$ws.Cells.Item(4, 18).Value = 50.09852060580601
$ws.Cells.Item(4, 19).Value = 0.1288364868576413
$ws.Cells.Item(4, 20).Value = 0.106140571625568

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Nid2"
$ws.Cells.Item(5, 3).Value = "Col13a1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 65.67978099999999
$ws.Cells.Item(5, 8).Value = 197.039343
$ws.Cells.Item(5, 9).Value = 0.6451826121651902
$ws.Cells.Item(5, 10).Value = 0.6640841926751532
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3015963333333334
$ws.Cells.Item(5, 14).Value = 0.9047890000000001
$ws.Cells.Item(5, 15).Value = 0.4737402899861982
$ws.Cells.Item(5, 16).Value = 0.5687660179545673
$ws.Cells.Item(5, 17).Value = 19.80878112373633
$ws.Cells.Item(5, 18).Value = 178.279030113627
$ws.Cells.Item(5, 19).Value = 0.30564899778119
$ws.Cells.Item(5, 20).Value = 0.3777085218544205

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Nid2"
$ws.Cells.Item(6, 3).Value = "Col13a1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 65.67978099999999
$ws.Cells.Item(6, 8).Value = 197.039343
$ws.Cells.Item(6, 9).Value = 0.6451826121651902
$ws.Cells.Item(6, 10).Value = 0.6640841926751532
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.01594066666666667
$ws.Cells.Item(6, 14).Value = 0.047822
$ws.Cells.Item(6, 15).Value = 0.02503921704145383
$ws.Cells.Item(6, 16).Value = 0.03006173650500096
$ws.Cells.Item(6, 17).Value = 1.046979495660667
$ws.Cells.Item(6, 18).Value = 9.422815460946
$ws.Cells.Item(6, 19).Value = 0.01615486745737633
$ws.Cells.Item(6, 20).Value = 0.01996352401733675

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Nid2"
$ws.Cells.Item(7, 3).Value = "Col13a1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 65.67978099999999
$ws.Cells.Item(7, 8).Value = 197.039343
$ws.Cells.Item(7, 9).Value = 0.6451826121651902
$ws.Cells.Item(7, 10).Value = 0.6640841926751532
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.319091
$ws.Cells.Item(7, 14).Value = 0.638182
$ws.Cells.Item(7, 15).Value = 0.5012204929723481
$ws.Cells.Item(7, 16).Value = 0.4011722455404317
$ws.Cells.Item(7, 17).Value = 20.957826999071
$ws.Cells.Item(7, 18).Value = 125.746961994426
$ws.Cells.Item(7, 19).Value = 0.3233787469266239
$ws.Cells.Item(7, 20).Value = 0.2664121468033959

# Row 8: M1 -> ECs
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Nid2"
$ws.Cells.Item(8, 3).Value = "Col13a1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6031273333333333
$ws.Cells.Item(8, 8).Value = 1.809382
$ws.Cells.Item(8, 9).Value = 0.005924612756979587
$ws.Cells.Item(8, 10).Value = 0.006098183065454873
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.3015963333333334
$ws.Cells.Item(8, 14).Value = 0.9047890000000001
$ws.Cells.Item(8, 15).Value = 0.4737402899861982
$ws.Cells.Item(8, 16).Value = 0.5687660179545673
$ws.Cells.Item(8, 17).Value = 0.1819009922664445
$ws.Cells.Item(8, 18).Value = 1.637108930398
$ws.Cells.Item(8, 19).Value = 0.002806727765547439
$ws.Cells.Item(8, 20).Value = 0.003468439298896744

# Row 9: M1 -> FAPs
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Nid2"
$ws.Cells.Item(9, 3).Value = "Col13a1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6031273333333333
$ws.Cells.Item(9, 8).Value = 1.809382
$ws.Cells.Item(9, 9).Value = 0.005924612756979587
$ws.Cells.Item(9, 10).Value = 0.006098183065454873
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.01594066666666667
$ws.Cells.Item(9, 14).Value = 0.047822
$ws.Cells.Item(9, 15).Value = 0.02503921704145383
$ws.Cells.Item(9, 16).Value = 0.03006173650500096
$ws.Cells.Item(9, 17).Value = 0.009614251778222224
$ws.Cells.Item(9, 18).Value = 0.08652826600400002
$ws.Cells.Item(9, 19).Value = 0.000148347664708578
$ws.Cells.Item(9, 20).Value = 0.0001833219724729634

# Row 10: M1 -> sCs
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Nid2"
$ws.Cells.Item(10, 3).Value = "Col13a1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.6031273333333333
$ws.Cells.Item(10, 8).Value = 1.809382
$ws.Cells.Item(10, 9).Value = 0.005924612756979587
$ws.Cells.Item(10, 10).Value = 0.006098183065454873
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.319091
$ws.Cells.Item(10, 14).Value = 0.638182
$ws.Cells.Item(10, 15).Value = 0.5012204929723481
$ws.Cells.Item(10, 16).Value = 0.4011722455404317
$ws.Cells.Item(10, 17).Value = 0.1924525039206667
$ws.Cells.Item(10, 18).Value = 1.154715023524
$ws.Cells.Item(10, 19).Value = 0.00296953732672357
$ws.Cells.Item(10, 20).Value = 0.002446421794085165

# Row 11: M2 -> ECs
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Nid2"
$ws.Cells.Item(11, 3).Value = "Col13a1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.6575653333333333
$ws.Cells.Item(11, 8).Value = 1.972696
$ws.Cells.Item(11, 9).Value = 0.006459365621655683
$ws.Cells.Item(11, 10).Value = 0.006648602307578259
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.3015963333333334
$ws.Cells.Item(11, 14).Value = 0.9047890000000001
$ws.Cells.Item(11, 15).Value = 0.4737402899861982
$ws.Cells.Item(11, 16).Value = 0.5687660179545673
$ws.Cells.Item(11, 17).Value = 0.1983192934604445
$ws.Cells.Item(11, 18).Value = 1.784873641144
$ws.Cells.Item(11, 19).Value = 0.003060061742730042
$ws.Cells.Item(11, 20).Value = 0.003781499059444833

# Row 12: M2 -> FAPs
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Nid2"
$ws.Cells.Item(12, 3).Value = "Col13a1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.6575653333333333
$ws.Cells.Item(12, 8).Value = 1.972696
$ws.Cells.Item(12, 9).Value = 0.006459365621655683
$ws.Cells.Item(12, 10).Value = 0.006648602307578259
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.01594066666666667
$ws.Cells.Item(12, 14).Value = 0.047822
$ws.Cells.Item(12, 15).Value = 0.02503921704145383
$ws.Cells.Item(12, 16).Value = 0.03006173650500096
$ws.Cells.Item(12, 17).Value = 0.01048202979022222
$ws.Cells.Item(12, 18).Value = 0.09433826811200001
$ws.Cells.Item(12, 19).Value = 0.000161737457750742
$ws.Cells.Item(12, 20).Value = 0.000199868530696959

# Row 13: M2 -> sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Nid2"
$ws.Cells.Item(13, 3).Value = "Col13a1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.6575653333333333
$ws.Cells.Item(13, 8).Value = 1.972696
$ws.Cells.Item(13, 9).Value = 0.006459365621655683
$ws.Cells.Item(13, 10).Value = 0.006648602307578259
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.319091
$ws.Cells.Item(13, 14).Value = 0.638182
$ws.Cells.Item(13, 15).Value = 0.5012204929723481
$ws.Cells.Item(13, 16).Value = 0.4011722455404317
$ws.Cells.Item(13, 17).Value = 0.2098231797786667
$ws.Cells.Item(13, 18).Value = 1.258939078672
$ws.Cells.Item(13, 19).Value = 0.003237566421174899
$ws.Cells.Item(13, 20).Value = 0.002667234717436466

# Row 14: sCs -> ECs
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Nid2"
$ws.Cells.Item(14, 3).Value = "Col13a1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 8.692511499999998
$ws.Cells.Item(14, 8).Value = 17.385023
$ws.Cells.Item(14, 9).Value = 0.08538788026479496
$ws.Cells.Item(14, 10).Value = 0.05859296315048091
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.3015963333333334
$ws.Cells.Item(14, 14).Value = 0.9047890000000001
$ws.Cells.Item(14, 15).Value = 0.4737402899861982
$ws.Cells.Item(14, 16).Value = 0.5687660179545673
$ws.Cells.Item(14, 17).Value = 2.621629595857833
$ws.Cells.Item(14, 18).Value = 15.729777575147
$ws.Cells.Item(14, 19).Value = 0.04045167915795073
$ws.Cells.Item(14, 20).Value = 0.03332568633125772

# Row 15: sCs -> FAPs
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Nid2"
$ws.Cells.Item(15, 3).Value = "Col13a1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 8.692511499999998
$ws.Cells.Item(15, 8).Value = 17.385023
$ws.Cells.Item(15, 9).Value = 0.08538788026479496
$ws.Cells.Item(15, 10).Value = 0.05859296315048091
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.01594066666666667
$ws.Cells.Item(15, 14).Value = 0.047822
$ws.Cells.Item(15, 15).Value = 0.02503921704145383
$ws.Cells.Item(15, 16).Value = 0.03006173650500096
$ws.Cells.Item(15, 17).Value = 0.1385644283176667
$ws.Cells.Item(15, 18).Value = 0.8313865699059999
$ws.Cells.Item(15, 19).Value = 0.002138045666659873
$ws.Cells.Item(15, 20).Value = 0.001761406219276988

# Row 16: sCs -> sCs
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Nid2"
$ws.Cells.Item(16, 3).Value = "Col13a1"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 8.692511499999998
$ws.Cells.Item(16, 8).Value = 17.385023
$ws.Cells.Item(16, 9).Value = 0.08538788026479496
$ws.Cells.Item(16, 10).Value = 0.05859296315048091
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.319091
$ws.Cells.Item(16, 14).Value = 0.638182
$ws.Cells.Item(16, 15).Value = 0.5012204929723481
$ws.Cells.Item(16, 16).Value = 0.4011722455404317
$ws.Cells.Item(16, 17).Value = 2.7737021870465
$ws.Cells.Item(16, 18).Value = 11.094808748186
$ws.Cells.Item(16, 19).Value = 0.04279815544018436
$ws.Cells.Item(16, 20).Value = 0.02350587059994619

Write-Output ("Used range after edit: " + $ws.UsedRange.Address())
